# Apply the "Added multiple unit support (basic)" change.
$wb = $excel.ActiveWorkbook

# 1. Rename worksheet "TEST_1.1" -> "Test_1.1"
$ws1 = $wb.ActiveSheet
[void]$ws1.Activate()
$ws1.Name = "Test_1.1"

# 2. Update D12:D15 on Test_1.1 to all reference the same unit ("Unit2/ENV", same as D11)
$unitValue = $ws1.Range("D11").Value2
$ws1.Range("D12").Value = $unitValue
$ws1.Range("D13").Value = $unitValue
$ws1.Range("D14").Value = $unitValue
$ws1.Range("D15").Value = $unitValue

# 3. Update the active selection on Test_1.1 to D11:D15
[void]$ws1.Range("D11:D15").Select()
